# Add the new "SOURCE_EXCLUDED_COLUMNS" key to the ETL example sheet.
#
# The original layout has columns A:I (SOURCE_* fields, ending with the
# recently-added SOURCE_TABLE_PRIMARY_KEY in I) immediately followed by the
# TARGET_*/etc. columns starting at J. This change inserts a brand-new
# column at J for SOURCE_EXCLUDED_COLUMNS, which pushes every column from
# the old J onward one slot to the right (J->K, K->L, ... Y->Z) along with
# their column-width formatting. Only the new header cell (J1) gets a
# value; the data rows (2-5) are left blank in the new column, matching
# every other "extra key" style column that has no per-row data yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting a whole column shifts existing data/formats right, same as
# right-clicking column J -> Insert in the Excel UI.
$ws.Columns("J:J").Insert()

# New column picks up the same width as its neighbors (H:I, width 38.5)
# instead of the sheet default.
$ws.Columns("J:J").ColumnWidth = $ws.Columns("I:I").ColumnWidth

# Header for the newly inserted column.
$ws.Range("J1").Value = "SOURCE_EXCLUDED_COLUMNS"

# Move the active selection to the new header cell.
$ws.Range("J1").Select() | Out-Null
